$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 1863.3334
$ws.Range("J55").Value = 5216.75
$ws.Range("L55").Value = 5216.75
$ws.Range("N55").Value = -5644.75
$ws.Range("H106").Value = 5999
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 5999
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 5999
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -7261
$ws.Range("H138").Value = 2072
$ws.Range("J138").Value = 3499.8
$ws.Range("L138").Value = 10499.4
$ws.Range("N138").Value = -20779.4
$ws.Range("H141").Value = 3346.9285
$ws.Range("I141").Value = 3346.9285
$ws.Range("K141").Value = 10040.7855
$ws.Range("M141").Value = -4860.7855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 20000
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8534
$ws.Range("I20").Value = 8534
$ws.Range("K20").Value = 8534
$ws.Range("M20").Value = -8287
$ws.Range("H21").Value = 60144
$ws.Range("J21").Value = 60144
$ws.Range("L21").Value = 60144
$ws.Range("N21").Value = -60616
$ws.Range("H28").Value = 60487.668
$ws.Range("J28").Value = 60487.668
$ws.Range("L28").Value = 60487.668
$ws.Range("N28").Value = -61075.668
$ws.Range("H86").Value = 3115.8667
$ws.Range("I86").Value = 2766
$ws.Range("J86").Value = 3515.7144
$ws.Range("K86").Value = 2766
$ws.Range("L86").Value = 3515.7144
$ws.Range("M86").Value = -1643
$ws.Range("N86").Value = -5761.7144
$ws.Range("H89").Value = 3115.8667
$ws.Range("I89").Value = 2766
$ws.Range("J89").Value = 3515.7144
$ws.Range("K89").Value = 13830
$ws.Range("L89").Value = 17578.572
$ws.Range("M89").Value = -8214
$ws.Range("N89").Value = -28810.572
$ws.Range("H111").Value = 104990
$ws.Range("J111").Value = 104990
$ws.Range("L111").Value = 104990
$ws.Range("N111").Value = -113170
$ws.Range("H112").Value = 79995
$ws.Range("J112").Value = 79995
$ws.Range("L112").Value = 79995
$ws.Range("N112").Value = -82949
$ws.Range("H117").Value = 49994.5
$ws.Range("J117").Value = 49994.5
$ws.Range("L117").Value = 49994.5
$ws.Range("N117").Value = -59172.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 48818
$ws.Range("I51").Value = 18199.6
$ws.Range("K51").Value = 18199.6
$ws.Range("M51").Value = -17463.6
$ws.Range("H61").Value = 48818
$ws.Range("I61").Value = 18199.6
$ws.Range("K61").Value = 18199.6
$ws.Range("M61").Value = -17851.6
$ws.Range("H68").Value = 77891.14
$ws.Range("J68").Value = 106996.664
$ws.Range("L68").Value = 106996.664
$ws.Range("N68").Value = -108494.664
$ws.Range("H71").Value = 77891.14
$ws.Range("J71").Value = 106996.664
$ws.Range("L71").Value = 320989.992
$ws.Range("N71").Value = -328477.992
$ws.Range("H75").Value = 113000
$ws.Range("J75").Value = 113000
$ws.Range("L75").Value = 113000
$ws.Range("N75").Value = -114996
$ws.Range("H78").Value = 113000
$ws.Range("J78").Value = 113000
$ws.Range("L78").Value = 339000
$ws.Range("N78").Value = -348984
$ws.Range("H107").Value = 776.34784
$ws.Range("I107").Value = 572.3684
$ws.Range("J107").Value = 1745.25
$ws.Range("K107").Value = 572.3684
$ws.Range("L107").Value = 1745.25
$ws.Range("M107").Value = 1347.6316
$ws.Range("N107").Value = -5585.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 872.75
$ws.Range("I14").Value = 872.75
$ws.Range("K14").Value = 2618.25
$ws.Range("M14").Value = -2445.25
$ws.Range("H70").Value = 500000350
$ws.Range("I70").Value = 500000350
$ws.Range("K70").Value = 1500001050
$ws.Range("M70").Value = -1500000735
$ws.Range("H73").Value = 500000350
$ws.Range("I73").Value = 500000350
$ws.Range("K73").Value = 1500001050
$ws.Range("M73").Value = -1499999958
$ws.Range("H140").Value = 2237.9
$ws.Range("I140").Value = 2237.9
$ws.Range("K140").Value = 6713.700000000001
$ws.Range("M140").Value = -1533.700000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3792.5
$ws.Range("I3").Value = 6752
$ws.Range("J3").Value = 833
$ws.Range("K3").Value = 6752
$ws.Range("L3").Value = 833
$ws.Range("M3").Value = -6636
$ws.Range("N3").Value = -1065
$ws.Range("H13").Value = 9129.286
$ws.Range("I13").Value = 10640.833
$ws.Range("J13").Value = 60
$ws.Range("K13").Value = 10640.833
$ws.Range("L13").Value = 60
$ws.Range("M13").Value = -10501.833
$ws.Range("N13").Value = -338
$ws.Range("H69").Value = 70018
$ws.Range("J69").Value = 70018
$ws.Range("L69").Value = 70018
$ws.Range("N69").Value = -71516
$ws.Range("H70").Value = 4936
$ws.Range("I70").Value = 4936
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4936
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4666
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 70018
$ws.Range("J72").Value = 70018
$ws.Range("L72").Value = 210054
$ws.Range("N72").Value = -217542
$ws.Range("H73").Value = 4936
$ws.Range("I73").Value = 4936
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4936
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -4000
$ws.Range("N73").Value = -4000
$ws.Range("H136").Value = 18154.2
$ws.Range("J136").Value = 18154.2
$ws.Range("L136").Value = 54462.60000000001
$ws.Range("N136").Value = -59562.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 39764.555
$ws.Range("I7").Value = 2090.611
$ws.Range("J7").Value = 115112.445
$ws.Range("K7").Value = 2090.611
$ws.Range("L7").Value = 115112.445
$ws.Range("M7").Value = -1978.611
$ws.Range("N7").Value = -115336.445
$ws.Range("H22").Value = 4000.3333
$ws.Range("I22").Value = 4000.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4000.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3705.3333
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 4000.3333
$ws.Range("I27").Value = 4000.3333
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 4000.3333
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3893.3333
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 4068.4119
$ws.Range("I40").Value = 2728.3333
$ws.Range("K40").Value = 2728.3333
$ws.Range("M40").Value = -2592.3333
$ws.Range("H45").Value = 3620.5
$ws.Range("I45").Value = 3620.5
$ws.Range("K45").Value = 3620.5
$ws.Range("M45").Value = -3213.5
$ws.Range("H126").Value = 39764.555
$ws.Range("I126").Value = 2090.611
$ws.Range("J126").Value = 115112.445
$ws.Range("K126").Value = 6271.833
$ws.Range("L126").Value = 345337.335
$ws.Range("M126").Value = -3801.833
$ws.Range("N126").Value = -350277.335
$ws.Range("H132").Value = 54457.414
$ws.Range("I132").Value = 28361.816
$ws.Range("J132").Value = 385001.66
$ws.Range("K132").Value = 85085.448
$ws.Range("L132").Value = 1155004.98
$ws.Range("M132").Value = -82555.448
$ws.Range("N132").Value = -1160064.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1369.2858
$ws.Range("I113").Value = 1549.25
$ws.Range("K113").Value = 4647.75
$ws.Range("M113").Value = -2477.75
$ws.Range("H122").Value = 3571.8206
$ws.Range("I122").Value = 2492.5186
$ws.Range("J122").Value = 6000.25
$ws.Range("K122").Value = 7477.5558
$ws.Range("L122").Value = 18000.75
$ws.Range("M122").Value = -5027.5558
$ws.Range("N122").Value = -22900.75
$ws.Range("H127").Value = 50000
$ws.Range("I127").Value = 50000
$ws.Range("K127").Value = 50000
$ws.Range("M127").Value = -45040

Write-Host "Edits applied successfully"
